$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-31 Thursday" "2025-08-01 Friday"

Replace-Text "640÷7=" "679÷5="
Replace-Text "297÷8=" "562÷9="
Replace-Text "399÷9=" "798÷6="
Replace-Text "290÷5=" "388÷8="
Replace-Text "732÷9=" "973÷4="

Replace-Text "542÷6=" "702÷7="
Replace-Text "254÷3=" "328÷4="
Replace-Text "912÷5=" "195÷9="
Replace-Text "103÷8=" "625÷5="
Replace-Text "132÷3=" "650÷7="

Replace-Text "654÷9=" "657÷9="
Replace-Text "755÷2=" "674÷3="
Replace-Text "638÷2=" "921÷6="
Replace-Text "594÷7=" "640÷7="
Replace-Text "885÷7=" "351÷3="

Replace-Text "171÷2=" "903÷2="
Replace-Text "335÷3=" "697÷9="
Replace-Text "814÷9=" "396÷2="
Replace-Text "769÷2=" "932÷2="
Replace-Text "512÷7=" "285÷5="

Replace-Text "537÷3=" "178÷8="
Replace-Text "216÷9=" "372÷8="
Replace-Text "634÷7=" "905÷5="
Replace-Text "700÷8=" "187÷7="
Replace-Text "737÷7=" "119÷6="

"Done"
